$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 1.090291
$ws.Range("N2").Value = 3.270873
$ws.Range("O2").Value = 0.02878663098006733
$ws.Range("P2").Value = 0.02878663098006733
$ws.Range("Q2").Value = 0.01582593729533333
$ws.Range("R2").Value = 0.142433435658
$ws.Range("S2").Value = 0.02878663098006733
$ws.Range("T2").Value = 0.02878663098006733

# Row 3 updates
$ws.Range("O3").Value = 0.7239320554917256
$ws.Range("P3").Value = 0.7239320554917257
$ws.Range("S3").Value = 0.7239320554917256
$ws.Range("T3").Value = 0.7239320554917257

# Row 4 updates
$ws.Range("O4").Value = 0.2472813135282071
$ws.Range("P4").Value = 0.2472813135282071
$ws.Range("S4").Value = 0.2472813135282071
$ws.Range("T4").Value = 0.2472813135282071
